# Apply odds updates to Sheet1 for the 2026-01-02 Betfair Back/Lay workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 4.8
$ws.Cells.Item(2, 10).Value = 4.1
$ws.Cells.Item(2, 11).Value = 4.2
$ws.Cells.Item(2, 14).Value = 4.8
$ws.Cells.Item(2, 16).Value = 2.24
$ws.Cells.Item(2, 17).Value = 1.74
$ws.Cells.Item(2, 18).Value = 1.5
$ws.Cells.Item(2, 19).Value = 2.9
$ws.Cells.Item(2, 20).Value = 1.73
$ws.Cells.Item(2, 21).Value = 2.3
$ws.Cells.Item(2, 25).Value = 20
$ws.Cells.Item(2, 40).Value = 9.800000000000001
$ws.Cells.Item(3, 7).Value = 4.4
$ws.Cells.Item(3, 12).Value = 1.32
$ws.Cells.Item(3, 17).Value = 1.7
$ws.Cells.Item(3, 20).Value = 1.65
$ws.Cells.Item(3, 23).Value = 1.3
$ws.Cells.Item(3, 25).Value = 12
$ws.Cells.Item(4, 6).Value = 1.62
$ws.Cells.Item(4, 8).Value = 6.6
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 11).Value = 4.4
$ws.Cells.Item(4, 12).Value = 1.31
$ws.Cells.Item(4, 15).Value = 1.3
$ws.Cells.Item(4, 16).Value = 1.95
$ws.Cells.Item(4, 17).Value = 1.88
$ws.Cells.Item(4, 18).Value = 1.35
$ws.Cells.Item(4, 20).Value = 1.97
$ws.Cells.Item(4, 21).Value = 1.91
$ws.Cells.Item(4, 22).Value = 1.17
$ws.Cells.Item(4, 24).Value = 16.5
$ws.Cells.Item(4, 25).Value = 23
$ws.Cells.Item(4, 27).Value = 260
$ws.Cells.Item(4, 32).Value = 11
$ws.Cells.Item(4, 34).Value = 26
$ws.Cells.Item(4, 36).Value = 15
$ws.Cells.Item(4, 39).Value = 180
$ws.Cells.Item(4, 40).Value = 9.4
$ws.Cells.Item(5, 6).Value = 2.5
$ws.Cells.Item(5, 7).Value = 2.68
$ws.Cells.Item(5, 8).Value = 3.45
$ws.Cells.Item(5, 9).Value = 3.9
$ws.Cells.Item(5, 10).Value = 2.82
$ws.Cells.Item(5, 17).Value = 2.96
$ws.Cells.Item(5, 19).Value = 7.2
$ws.Cells.Item(5, 20).Value = 2.28
$ws.Cells.Item(5, 21).Value = 1.65
$ws.Cells.Item(5, 23).Value = 1.59
$ws.Cells.Item(5, 24).Value = 7.8
$ws.Cells.Item(5, 25).Value = 9.6
$ws.Cells.Item(5, 26).Value = 22
$ws.Cells.Item(5, 27).Value = 85
$ws.Cells.Item(5, 29).Value = 6.8
$ws.Cells.Item(5, 31).Value = 1000
$ws.Cells.Item(5, 33).Value = 14
$ws.Cells.Item(5, 34).Value = 30
$ws.Cells.Item(5, 37).Value = 46
$ws.Cells.Item(5, 38).Value = 95
$ws.Cells.Item(5, 39).Value = 1000
$ws.Cells.Item(5, 40).Value = 60
$ws.Cells.Item(6, 6).Value = 2.4
$ws.Cells.Item(6, 7).Value = 2.68
$ws.Cells.Item(6, 8).Value = 3.25
$ws.Cells.Item(6, 9).Value = 3.95
$ws.Cells.Item(6, 11).Value = 3.45
$ws.Cells.Item(6, 13).Value = 1.07
$ws.Cells.Item(6, 14).Value = 2.54
$ws.Cells.Item(6, 15).Value = 1.07
$ws.Cells.Item(6, 17).Value = 2.3
$ws.Cells.Item(6, 18).Value = 1.18
$ws.Cells.Item(6, 20).Value = 1.05
$ws.Cells.Item(6, 22).Value = 1.37
$ws.Cells.Item(6, 23).Value = 1.6
$ws.Cells.Item(7, 8).Value = 1.58
$ws.Cells.Item(7, 9).Value = 1.65
$ws.Cells.Item(7, 14).Value = 6
$ws.Cells.Item(7, 16).Value = 2.72
$ws.Cells.Item(7, 18).Value = 1.69
$ws.Cells.Item(7, 20).Value = 1.6
$ws.Cells.Item(7, 21).Value = 2.36
$ws.Cells.Item(7, 26).Value = 13.5
$ws.Cells.Item(7, 34).Value = 990
$ws.Cells.Item(8, 8).Value = 1.5
$ws.Cells.Item(8, 9).Value = 1.51
$ws.Cells.Item(8, 10).Value = 4.6
$ws.Cells.Item(8, 14).Value = 3.55
$ws.Cells.Item(8, 15).Value = 1.36
$ws.Cells.Item(8, 16).Value = 1.86
$ws.Cells.Item(8, 17).Value = 2.08
$ws.Cells.Item(8, 18).Value = 1.33
$ws.Cells.Item(8, 19).Value = 3.75
$ws.Cells.Item(8, 20).Value = 2.22
$ws.Cells.Item(8, 21).Value = 1.76
$ws.Cells.Item(8, 25).Value = 7.2
$ws.Cells.Item(8, 28).Value = 22
$ws.Cells.Item(8, 31).Value = 17.5
$ws.Cells.Item(8, 38).Value = 160
$ws.Cells.Item(8, 39).Value = 210
$ws.Cells.Item(8, 40).Value = 240
$ws.Cells.Item(8, 41).Value = 9.6
$ws.Cells.Item(9, 11).Value = 3.85
$ws.Cells.Item(9, 20).Value = 2.22
$ws.Cells.Item(9, 32).Value = 8.800000000000001
$ws.Cells.Item(9, 38).Value = 980
$ws.Cells.Item(10, 6).Value = 2.74
$ws.Cells.Item(10, 7).Value = 2.76
$ws.Cells.Item(10, 8).Value = 2.8
$ws.Cells.Item(10, 9).Value = 2.82
$ws.Cells.Item(10, 10).Value = 3.5
$ws.Cells.Item(10, 16).Value = 1.94
$ws.Cells.Item(10, 21).Value = 2.2
$ws.Cells.Item(10, 22).Value = 1.54
$ws.Cells.Item(10, 26).Value = 18.5
$ws.Cells.Item(10, 37).Value = 30
$ws.Cells.Item(10, 40).Value = 26
$ws.Cells.Item(11, 6).Value = 6.8
$ws.Cells.Item(11, 7).Value = 7
$ws.Cells.Item(11, 8).Value = 1.59
$ws.Cells.Item(11, 9).Value = 1.6
$ws.Cells.Item(11, 10).Value = 4.4
$ws.Cells.Item(11, 11).Value = 4.5
$ws.Cells.Item(11, 16).Value = 1.96
$ws.Cells.Item(11, 17).Value = 1.99
$ws.Cells.Item(11, 21).Value = 1.89
$ws.Cells.Item(11, 22).Value = 2.66
$ws.Cells.Item(11, 23).Value = 1.16
$ws.Cells.Item(11, 24).Value = 14.5
$ws.Cells.Item(11, 25).Value = 7.6
$ws.Cells.Item(11, 30).Value = 9.800000000000001
$ws.Cells.Item(12, 14).Value = 2.28
$ws.Cells.Item(12, 16).Value = 1.41
$ws.Cells.Item(12, 19).Value = 7.6
$ws.Cells.Item(12, 23).Value = 1.82
$ws.Cells.Item(12, 24).Value = 6.8
$ws.Cells.Item(12, 28).Value = 6
$ws.Cells.Item(13, 6).Value = 1.87
$ws.Cells.Item(13, 7).Value = 1.91
$ws.Cells.Item(13, 10).Value = 3.6
$ws.Cells.Item(13, 11).Value = 3.7
$ws.Cells.Item(13, 17).Value = 2.16
$ws.Cells.Item(13, 23).Value = 2.08
$ws.Cells.Item(13, 27).Value = 140
$ws.Cells.Item(13, 30).Value = 21
$ws.Cells.Item(13, 41).Value = 110
